# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Reorder "Guyana" / "Islas Caimanes" in the country list (rows 150-151) ---
# Previously: row150 = Islas Caimanes, row151 = Guyana
# Now:        row150 = Guyana,         row151 = Islas Caimanes
$ws.Range("A150").Value = "Guyana"
$ws.Range("A151").Value = "Islas Caimanes"

# --- Update case numbers (Casos totales, Nuevos casos, Casos activos, Recuperados,
#     Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 700234
$ws.Range("C4").Value = 22664
$ws.Range("D4").Value = 59158
$ws.Range("E4").Value = 604154
$ws.Range("G4").Value = 2305
$ws.Range("H4").Value = 36922

# Row 5: Espana
$ws.Range("B5").Value = 188167
$ws.Range("C5").Value = 3219
$ws.Range("E5").Value = 93739
$ws.Range("G5").Value = 316
$ws.Range("H5").Value = 19631

# Row 10: China
$ws.Range("C10").Value = 325

# Row 20: Austria
$ws.Range("B20").Value = 14595
$ws.Range("C20").Value = 119
$ws.Range("E20").Value = 4460

# Row 150: now Guyana (new data)
$ws.Range("B150").Value = 63
$ws.Range("C150").Value = 8
$ws.Range("D150").Value = 9
$ws.Range("E150").Value = 48
$ws.Range("F150").Value = 4
$ws.Range("H150").Value = 6

# Row 151: now Islas Caimanes (same data previously held by the old row150)
$ws.Range("B151").Value = 61
$ws.Range("C151").Value = 1
$ws.Range("D151").Value = 7
$ws.Range("E151").Value = 53
$ws.Range("F151").Value = 3
$ws.Range("H151").Value = 1
